$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "Planning"

# Add "Unnamed" header labels for columns G:J, matching the style of existing headers (A1:F1)
$ws.Range("A1").Copy() | Out-Null
$ws.Range("G1:J1").PasteSpecial(-4122) | Out-Null # xlPasteFormats

$ws.Range("G1").Value = "Unnamed: 6"
$ws.Range("H1").Value = "Unnamed: 7"
$ws.Range("I1").Value = "Unnamed: 8"
$ws.Range("J1").Value = "Unnamed: 9"

# Project title text
$ws.Range("I2").Value = "Project title:"
$ws.Range("J2").Value = "My Gantt diagram"

$wb.Save()
